# Natmi following Dr Hou advice
# Expand the Ncam1-Robo1 sending/target cluster matrix from 2x2 (FAPs/sCs)
# to a full 3x3 matrix that also includes the "ECs" cluster.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 2: ECs -> ECs (was FAPs -> FAPs) ---
$ws.Range("A2").Value = "ECs"
$ws.Range("D2").Value = "ECs"
$ws.Range("E2").Value = 2
$ws.Range("F2").Value = 0.6666666666666666
$ws.Range("G2").Value = 0.9404873333333333
$ws.Range("H2").Value = 2.821462
$ws.Range("I2").Value = 0.02000383747045655
$ws.Range("J2").Value = 0.02000383747045654
$ws.Range("K2").Value = 2
$ws.Range("L2").Value = 0.6666666666666666
$ws.Range("M2").Value = 0.1325553333333333
$ws.Range("N2").Value = 0.397666
$ws.Range("O2").Value = 0.00533964316398423
$ws.Range("P2").Value = 0.00533964316398423
$ws.Range("Q2").Value = 0.1246666119657778
$ws.Range("R2").Value = 1.121999507692
$ws.Range("S2").Value = 0.0001068133540025749
$ws.Range("T2").Value = 0.0001068133540025749

# --- Row 3: ECs -> FAPs (was FAPs -> sCs) ---
$ws.Range("A3").Value = "ECs"
$ws.Range("D3").Value = "FAPs"
$ws.Range("E3").Value = 2
$ws.Range("F3").Value = 0.6666666666666666
$ws.Range("G3").Value = 0.9404873333333333
$ws.Range("H3").Value = 2.821462
$ws.Range("I3").Value = 0.02000383747045655
$ws.Range("J3").Value = 0.02000383747045654
$ws.Range("K3").Value = 3
$ws.Range("L3").Value = 1
$ws.Range("M3").Value = 17.178266
$ws.Range("N3").Value = 51.534798
$ws.Range("O3").Value = 0.6919812904497951
$ws.Range("P3").Value = 0.691981290449795
$ws.Range("Q3").Value = 16.15594158163067
$ws.Range("R3").Value = 145.403474234676
$ws.Range("S3").Value = 0.01384228126675448
$ws.Range("T3").Value = 0.01384228126675448

# --- Row 4: ECs -> sCs (was sCs -> FAPs) ---
$ws.Range("A4").Value = "ECs"
$ws.Range("B4").Value = "Ncam1"
$ws.Range("C4").Value = "Robo1"
$ws.Range("D4").Value = "sCs"
$ws.Range("E4").Value = 2
$ws.Range("F4").Value = 0.6666666666666666
$ws.Range("G4").Value = 0.9404873333333333
$ws.Range("H4").Value = 2.821462
$ws.Range("I4").Value = 0.02000383747045655
$ws.Range("J4").Value = 0.02000383747045654
$ws.Range("K4").Value = 3
$ws.Range("L4").Value = 1
$ws.Range("M4").Value = 7.513933666666667
$ws.Range("N4").Value = 22.541801
$ws.Range("O4").Value = 0.3026790663862208
$ws.Range("P4").Value = 0.3026790663862208
$ws.Range("Q4").Value = 7.066759437006889
$ws.Range("R4").Value = 63.600834933062
$ws.Range("S4").Value = 0.006054742849699487
$ws.Range("T4").Value = 0.006054742849699486

# --- Row 5: FAPs -> ECs (new) ---
$ws.Range("A5").Value = "FAPs"
$ws.Range("B5").Value = "Ncam1"
$ws.Range("C5").Value = "Robo1"
$ws.Range("D5").Value = "ECs"
$ws.Range("E5").Value = 3
$ws.Range("F5").Value = 1
$ws.Range("G5").Value = 1.392600333333333
$ws.Range("H5").Value = 4.177801000000001
$ws.Range("I5").Value = 0.0296201232509638
$ws.Range("J5").Value = 0.0296201232509638
$ws.Range("K5").Value = 2
$ws.Range("L5").Value = 0.6666666666666666
$ws.Range("M5").Value = 0.1325553333333333
$ws.Range("N5").Value = 0.397666
$ws.Range("O5").Value = 0.00533964316398423
$ws.Range("P5").Value = 0.00533964316398423
$ws.Range("Q5").Value = 0.1845966013851111
$ws.Range("R5").Value = 1.661369412466
$ws.Range("S5").Value = 0.0001581608886333792
$ws.Range("T5").Value = 0.0001581608886333792

# --- Row 6: FAPs -> FAPs (new) ---
$ws.Range("A6").Value = "FAPs"
$ws.Range("B6").Value = "Ncam1"
$ws.Range("C6").Value = "Robo1"
$ws.Range("D6").Value = "FAPs"
$ws.Range("E6").Value = 3
$ws.Range("F6").Value = 1
$ws.Range("G6").Value = 1.392600333333333
$ws.Range("H6").Value = 4.177801000000001
$ws.Range("I6").Value = 0.0296201232509638
$ws.Range("J6").Value = 0.0296201232509638
$ws.Range("K6").Value = 3
$ws.Range("L6").Value = 1
$ws.Range("M6").Value = 17.178266
$ws.Range("N6").Value = 51.534798
$ws.Range("O6").Value = 0.6919812904497951
$ws.Range("P6").Value = 0.691981290449795
$ws.Range("Q6").Value = 23.92245895768867
$ws.Range("R6").Value = 215.302130619198
$ws.Range("S6").Value = 0.02049657111048391
$ws.Range("T6").Value = 0.0204965711104839

# --- Row 7: FAPs -> sCs (new) ---
$ws.Range("A7").Value = "FAPs"
$ws.Range("B7").Value = "Ncam1"
$ws.Range("C7").Value = "Robo1"
$ws.Range("D7").Value = "sCs"
$ws.Range("E7").Value = 3
$ws.Range("F7").Value = 1
$ws.Range("G7").Value = 1.392600333333333
$ws.Range("H7").Value = 4.177801000000001
$ws.Range("I7").Value = 0.0296201232509638
$ws.Range("J7").Value = 0.0296201232509638
$ws.Range("K7").Value = 3
$ws.Range("L7").Value = 1
$ws.Range("M7").Value = 7.513933666666667
$ws.Range("N7").Value = 22.541801
$ws.Range("O7").Value = 0.3026790663862208
$ws.Range("P7").Value = 0.3026790663862208
$ws.Range("Q7").Value = 10.46390652884456
$ws.Range("R7").Value = 94.175158759601
$ws.Range("S7").Value = 0.008965391251846514
$ws.Range("T7").Value = 0.008965391251846512

# --- Row 8: sCs -> ECs (new) ---
$ws.Range("A8").Value = "sCs"
$ws.Range("B8").Value = "Ncam1"
$ws.Range("C8").Value = "Robo1"
$ws.Range("D8").Value = "ECs"
$ws.Range("E8").Value = 3
$ws.Range("F8").Value = 1
$ws.Range("G8").Value = 44.682258
$ws.Range("H8").Value = 134.046774
$ws.Range("I8").Value = 0.9503760392785797
$ws.Range("J8").Value = 0.9503760392785796
$ws.Range("K8").Value = 2
$ws.Range("L8").Value = 0.6666666666666666
$ws.Range("M8").Value = 0.1325553333333333
$ws.Range("N8").Value = 0.397666
$ws.Range("O8").Value = 0.00533964316398423
$ws.Range("P8").Value = 0.00533964316398423
$ws.Range("Q8").Value = 5.922871603276
$ws.Range("R8").Value = 53.305844429484
$ws.Range("S8").Value = 0.005074668921348276
$ws.Range("T8").Value = 0.005074668921348276

# --- Row 9: sCs -> FAPs (was sCs -> FAPs in row 4 before, now a distinct new row) ---
$ws.Range("A9").Value = "sCs"
$ws.Range("B9").Value = "Ncam1"
$ws.Range("C9").Value = "Robo1"
$ws.Range("D9").Value = "FAPs"
$ws.Range("E9").Value = 3
$ws.Range("F9").Value = 1
$ws.Range("G9").Value = 44.682258
$ws.Range("H9").Value = 134.046774
$ws.Range("I9").Value = 0.9503760392785797
$ws.Range("J9").Value = 0.9503760392785796
$ws.Range("K9").Value = 3
$ws.Range("L9").Value = 1
$ws.Range("M9").Value = 17.178266
$ws.Range("N9").Value = 51.534798
$ws.Range("O9").Value = 0.6919812904497951
$ws.Range("P9").Value = 0.691981290449795
$ws.Range("Q9").Value = 767.5637134046279
$ws.Range("R9").Value = 6908.073420641652
$ws.Range("S9").Value = 0.6576424380725567
$ws.Range("T9").Value = 0.6576424380725565

# --- Row 10: sCs -> sCs (was sCs -> sCs in row 5 before, now re-derived) ---
$ws.Range("A10").Value = "sCs"
$ws.Range("B10").Value = "Ncam1"
$ws.Range("C10").Value = "Robo1"
$ws.Range("D10").Value = "sCs"
$ws.Range("E10").Value = 3
$ws.Range("F10").Value = 1
$ws.Range("G10").Value = 44.682258
$ws.Range("H10").Value = 134.046774
$ws.Range("I10").Value = 0.9503760392785797
$ws.Range("J10").Value = 0.9503760392785796
$ws.Range("K10").Value = 3
$ws.Range("L10").Value = 1
$ws.Range("M10").Value = 7.513933666666667
$ws.Range("N10").Value = 22.541801
$ws.Range("O10").Value = 0.3026790663862208
$ws.Range("P10").Value = 0.3026790663862208
$ws.Range("Q10").Value = 335.739522688886
$ws.Range("R10").Value = 3021.655704199974
$ws.Range("S10").Value = 0.2876589322846748
$ws.Range("T10").Value = 0.2876589322846748
